# docs/protocol.xlsx
#
# "don't allow you_close_it for C2S HTTP. you_close_it will be the default
#  behavior if client doesn't send gimme_boxes."
#
# The "you_close_it" row (row 16) on the "Protocol" sheet previously allowed
# you_close_it for both S2C HTTP and C2S HTTP transports, with a comment
# describing both contexts. Now it is S2C-only: the C2S HTTP column flips
# from "Y" to "N", and the trailing "In C2S context: ..." paragraph is
# dropped from the comment, leaving only the S2C explanation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

# Column B is "C2S HTTP" support for this message - no longer allowed.
$ws.Range("B16").Value = "N"

# Trim the comment down to just the S2C-context paragraph.
$ws.Range("F16").Value = "In S2C context: This is useful when server wants client to do active close on the TCP socket for this transport (usually, to avoid having TIME_WAIT sockets). If server sends you_close_it over an HTTP transport, client must try to close the HTTP connection. (usually with xhrObject.abort() or removing an iframe)"

# The row shrinks now that the comment is a single paragraph instead of two.
$ws.Rows.Item(16).RowHeight = 40.5

# The active selection moved from A7 to F7.
$ws.Activate()
$ws.Range("F7").Select() | Out-Null
